$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Foglio1")
$ws2 = $wb.Worksheets.Item("Foglio2")

# --- Foglio1: make room for the two new "ModifyStudentsAlreadyAddedExternally" rows ---
# One new row inside the Entrance block (pushes NullPointerExceptionTest from row 9 to row 10)
$ws1.Rows("9").Insert()
# One new row inside the DR block (pushes the blank separator / Dashboard block down by one)
$ws1.Rows("17").Insert()

# Restore the row height on the two freshly inserted rows to match their neighbours
$ws1.Rows("9").RowHeight = 15.75
$ws1.Rows("17").RowHeight = 15.75

# --- Fix up merged cells for the class labels (A column) ---
$ws1.Range("A11:A15").UnMerge()
$ws1.Range("A2:A10").Merge()
$ws1.Range("A12:A17").Merge()
# A17:A29 shifted automatically with the row inserts to A19:A31 - nothing else to do there

# --- New Entrance row 9: ModifyStudentsAlreadyAddedExternally ---
$ws1.Cells.Item(9,2).Value = "ModifyStudentsAlreadyAddedExternally"
$ws1.Cells.Item(9,3).Value = "FAILED/USELESS"
$ws1.Cells.Item(9,4).Value = "Bisogna controllare che non venga modificato uno studente inserito nella entrance modificando lo studente passato, inutile, lo studente non può essere modificato"

# --- Entrance exceptions previously FAILED are now FAILED/USELESS ---
$ws1.Cells.Item(4,3).Value = "FAILED/USELESS"
$ws1.Cells.Item(5,3).Value = "FAILED/USELESS"

# --- New DR row 17: ModifyStudentsAlreadyAddedExternally ---
$ws1.Cells.Item(17,2).Value = "ModifyStudentsAlreadyAddedExternally"
$ws1.Cells.Item(17,3).Value = "FAILED/USELESS"
$ws1.Cells.Item(17,4).Value = "Bisogna controllare che non venga modificato uno studente inserito nella entrance modificando il set usato per inserirlo"

# --- Dashboard section: tests moved from TO DO to PASSED, and a renamed test ---
$ws1.Cells.Item(23,3).Value = "PASSED"   # InsertAndRemoveFromEntranceTest

$ws1.Cells.Item(24,2).Value = "MoveFromEntranceTODRTest"
$ws1.Cells.Item(24,3).Value = "PASSED"
$ws1.Cells.Item(24,4).Value = "Usato per controllare che gli studenti vengano inseriti nella DR"

$ws1.Cells.Item(25,4).Value = "Molto probabilmente eccezione inutile"  # WrongColorExceptionTest comment

$ws1.Cells.Item(26,3).Value = "PASSED"   # CardTest
$ws1.Cells.Item(27,3).Value = "PASSED"   # MasterTest

# --- sheet view / selection ---
$ws1.Activate()
$ws1.Range("B28").Select()

# --- Foglio2: widen column A, change A4 text, update selection ---
$ws2.Columns("A").ColumnWidth = 27.45
$ws2.Cells.Item(4,1).Value = "FAILED/USELESS"
$ws2.Activate()
$ws2.Range("D4").Select()

$ws1.Activate()
